# "further cleaning to metadata"
#  - libraryProtocol value used by every data row (K2:K29) changes from
#    "E7760" to "E7420"
#  - the roboticLibraryPrep column (L2:L29) is re-expressed as a FALSE()
#    formula instead of a bare boolean literal
#  - the sheet's remembered selection moves from L2:L29 to K2:K29

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the libraryProtocol text in place (keeps it in the shared-string
# table rather than creating a disconnected duplicate).
$null = $ws.Cells.Replace("E7760", "E7420")

# Re-enter the roboticLibraryPrep flag for every data row as a formula.
# Writing cell-by-cell (rather than to the whole L2:L29 range at once)
# keeps each cell's own formula instead of Excel collapsing them into one
# shared formula spanning the range.
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 12).Formula = "=FALSE()"
}

# Match the new remembered selection.
$null = $ws.Range("K2:K29").Select()
